$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 5).Value = "SIN REGISTRO $n"
}
